$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B/C (plain text: coin names & links) ---
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("B21").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C21").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

# --- Column D/E (numeric-looking text: must force text storage) ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.52'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.95%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '30.32'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '11.43%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.172'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.16%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05730'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.73%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.602'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.97%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.071'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.16%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8563'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '4.97%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8802'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '6.02%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01025'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1,605.66%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1366'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.80%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07144'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '3.16%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02865'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-2.34%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09393'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.03%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001520'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.24%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04146'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.06%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006027'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.68%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.490'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.25%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.266'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.87%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3214'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '3.24%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.03248'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '4.51%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1300'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.70%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.509'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-6.10%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.47%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.96%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004493'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.39%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '0.01%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03785'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.55%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005680'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-6.08%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1070'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.69%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002200'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.31%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01002'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '21.79%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005089'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-5.62%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.04%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.08000'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-40.72%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002763'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-41.40%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.04%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.04%'
